$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; existing rows 12-15 shift down to 13-16.
$ws.Rows("12:12").Insert()

# Populate the new row 12 with the new data entry.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 44879
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100104
$ws.Range("H12").Value = "Frutos de pepita"
$ws.Range("I12").Value = 100104004
$ws.Range("J12").Value = "Níspero"
$ws.Range("K12").Value = "Californiana(o)"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 25
$ws.Range("N12").Value = 30000
$ws.Range("O12").Value = 30000
$ws.Range("P12").Value = 30000
$ws.Range("Q12").Value = "$/bandeja 10 kilos"
$ws.Range("R12").Value = "Provincia de Quillota"
$ws.Range("S12").Value = 3000
$ws.Range("T12").Value = 10

# Ensure the date cell keeps the date/time number format used by the rest of column D.
$ws.Range("D12").NumberFormat = $ws.Range("D13").NumberFormat
